$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the IEEE reference-guide URL (old http journals.* link -> new
#    https ieeeauthorcenter.* link) and give the URL its own Courier New,
#    9pt run, matching the target formatting.
# ---------------------------------------------------------------------------
$oldUrl = "http://journals.ieeeauthorcenter.ieee.org/wp-content/uploads/sites/7/IEEE-Reference-Guide.pdf"
$newUrl = "https://ieeeauthorcenter.ieee.org/wp-content/uploads/IEEE-Reference-Guide.pdf"

$rng = $d.Content
$found = $rng.Find.Execute("must conform to the IEEE standard (" + $oldUrl + "). ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "must conform to the IEEE standard (" + $newUrl + "). "

    $urlRng = $d.Content
    $urlRng.Start = $rng.Start
    $urlRng.End = $rng.End
    $foundUrl = $urlRng.Find.Execute($newUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundUrl) {
        $urlRng.Font.Name = "Courier New"
        $urlRng.Font.Size = 9
    }
}

# ---------------------------------------------------------------------------
# 2) Fix the "A. Author and B Authour" typo in the first reference so it
#    reads "A. Author and B. Author".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("A. Author and B Authour, ", $true, $false, $false, $false, $false, $true, 1, $false, "A. Author and B. Author, ", 2) | Out-Null

$d.Save()
